# KANSAS_2018.xlsx cleanup:
#  1. Rename header columns to snake_case English names.
#  2. Title-case the lowercase Spanish connector words ("de", "del", "la",
#     "las", "el", "los", "y") that appear inside the mx_state / mx_municipality
#     text values (columns A and B).
#  3. Remove the trailing metadata/footnote rows (964-968) and let the used
#     range / dimension shrink back down to A1:D962.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header rename -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A & B (rows 2-962) ------
$connectors = @("de", "del", "la", "las", "el", "los", "y")

for ($r = 2; $r -le 962; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val -ne "") {
            $words = $val -split " "
            $changed = $false
            for ($i = 0; $i -lt $words.Length; $i++) {
                $w = $words[$i]
                if ($connectors -contains $w.ToLower()) {
                    $words[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1)
                    $changed = $true
                }
            }
            if ($changed) {
                $cell.Value = [string]::Join(" ", $words)
            }
        }
    }
}

# --- 3. Drop the trailing metadata rows (964-968) ----------------------
$ws.Range("A964:A968").EntireRow.Delete() | Out-Null
